# Fruta / hortaliza, semanal
# Insert a new weekly record at row 112 ("Americana O Klondike", Primera,
# fecha 2021-11-08) on the "Sandia" sheet, pushing the existing rows
# 112-134 down to 113-135.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 112; this shifts rows
# 112..134 down to 113..135 and extends the sheet dimension to R135.
$ws.Rows(112).Insert()

$ws.Cells.Item(112, 1).Value = 4
$ws.Cells.Item(112, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(112, 3).Value = "Los Lagos"
$ws.Cells.Item(112, 4).Value = 44508
$ws.Cells.Item(112, 5).Value = 10
$ws.Cells.Item(112, 6).Value = 100112028
$ws.Cells.Item(112, 7).Value = "Sandia"
$ws.Cells.Item(112, 8).Value = "Americana O Klondike"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 350
$ws.Cells.Item(112, 11).Value = 1000
$ws.Cells.Item(112, 12).Value = 1000
$ws.Cells.Item(112, 13).Value = 1000
$ws.Cells.Item(112, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(112, 15).Value = "Perú"
$ws.Cells.Item(112, 16).Value = 1000
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = "Hortaliza"
